# Strategy.xlsx - "FINAL TOUCHES ADDED - FIXED FILE PATHS OF ALL FILES"
# Adds three new games (Hexa Sort, Words of Wonders, Mergest Kingdom) to the
# "Strategy" sheet (rows 12-14: NAME, LINK, ABOUT, IMAGE SRC columns) and
# switches the active tab back to "Strategy" (away from "All").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Strategy")

# --- Row 12: Hexa Sort -------------------------------------------------
# Cell write order matters for shared-string allocation order: D, A, B, C.
$ws.Range("D12").Value = "`"https://images.crazygames.com/hexa-sort_16x9/20240…at%2Ccompress&q=65&cs=strip&ch=DPR&w=710&fit=crop`""
$ws.Range("A12").Value = "Hexa Sort"
$ws.Range("B12").Value = "https://www.crazygames.com/game/hexa-sort"
$ws.Range("C12").Value = "Hexa Sort is a captivating puzzle game that combines strategic matching and merging challenges. Engage in brain-teasing puzzles where you shuffle and organize hexagon tiles to achieve satisfying color matches. With smooth 3D graphics, vibrant colors, and relaxing gameplay, Hexa Sort offers a perfect balance of excitement and calm, making it ideal for both challenge seekers and those looking to unwind."

# --- Row 13: Words of Wonders ------------------------------------------
# Cell write order: D, B, A, C. Column A also picks up the bold "name"
# style used elsewhere in the sheet (style index 1, copied from A2's
# formatting so no new font/style entries are created).
$ws.Range("D13").Value = "`"https://images.crazygames.com/words-of-wonders_16x…at%2Ccompress&q=65&cs=strip&ch=DPR&w=336&fit=crop`""
$ws.Range("B13").Value = "https://www.crazygames.com/game/words-of-wonders"
$ws.Range("A13").Value = "Words of Wonder"
$ws.Range("C13").Value = "Words of Wonders is a crossword puzzle where you connect letters to reveal hidden words. Test your vocabulary and find all the words to progress and reach more challenging levels. Will you become a master wordsmith in this delightful word game?"

$ws.Range("A2").Copy()
$ws.Range("A13").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 14: Mergest Kingdom --------------------------------------------
# Cell write order: D, B, A, C.
$ws.Range("D14").Value = "`"https://images.crazygames.com/mergest-kingdom_16x9…at%2Ccompress&q=65&cs=strip&ch=DPR&w=336&fit=crop`""
$ws.Range("B14").Value = "https://www.crazygames.com/game/mergest-kingdom"
$ws.Range("A14").Value = "Mergest Kingdom"
$ws.Range("C14").Value = "Mergest Kingdom is a charming merge puzzle game that lets you build your own kingdom by matching various resources in groups of 3. Start building your fairytale kingdom by completing quests and matching hundreds of objects!"

# --- View state: switch focus back to the Strategy tab -----------------
$ws.Activate()
$ws.Range("A15").Select()
